# Auto-generated edit script: refresh cryptos price/volume table with latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Note: the source data stores Price (column D) as literal text, even when a price
# looks like a plain number (e.g. "0.9999", "112.20", "69.50"). Excels COM Value
# setter auto-converts number-looking strings to real numbers (dropping trailing
# zeros / "thousand dot" formatting), so for those cells we briefly force the Text
# number format, assign the literal string, then restore the Normal style so the
# cell keeps its original (unstyled) appearance with the exact text preserved.

$ws.Range('D2').Value = '30.580.13'
$ws.Range('E2').Value = '  -0.06%  '

$ws.Range('D3').Value = '1.920.96'
$ws.Range('E3').Value = '  -0.12%  '

$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '0.9999'
$ws.Range('D4').Style = "Normal"
$ws.Range('E4').Value = '  -0.06%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '245.42'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -0.78%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.9997'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -0.06%  '

$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.4826'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  +1.71%  '

$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.2900'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  -0.05%  '

$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.06810'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -0.21%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '112.20'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +6.45%  '

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '19.45'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +5.72%  '

$ws.Range('D12').Value = '1.913.39'
$ws.Range('E12').Value = '  -0.58%  '

$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '5.494'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +2.56%  '

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.07571'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -1.68%  '

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.6736'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +0.55%  '

$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '295.08'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +1.47%  '

$ws.Range('D17').Value = '30.569.53'
$ws.Range('E17').Value = '  -0.13%  '

$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.000007672'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  +0.65%  '

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '13.04'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +0.66%  '

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '0.9997'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -0.07%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '5.516'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -0.76%  '

$ws.Range('D22').Value = '2.160.47'
$ws.Range('E22').Value = '  -0.73%  '

$ws.Range('E23').Value = '  -0.16%  '

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '6.448'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -0.24%  '

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '9.496'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -0.22%  '

$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '167.11'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -0.34%  '

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '20.36'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -2.76%  '

$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '2.093'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -1.42%  '

$ws.Range('E29').Value = '  -0.62%  '

$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '1.442'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +2.72%  '

$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '4.138'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -0.97%  '

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '4.061'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +0.28%  '

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.04985'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -0.58%  '

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.7345'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +0.09%  '

$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.138'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -0.68%  '

$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '2.714'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -0.30%  '

$ws.Range('E37').Value = '  -2.08%  '

$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '2.683'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -0.12%  '

$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '2.025'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -0.73%  '

$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '109.50'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -2.01%  '

$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.4438'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +0.64%  '

$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.8693'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -0.59%  '

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '5.867'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -0.59%  '

$ws.Range('E44').Value = '  -0.02%  '

$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '69.50'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +2.53%  '

$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '7.256'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -0.47%  '

$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '49.01'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +1.60%  '

$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '9.236'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -1.12%  '

$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.1229'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -1.12%  '

$ws.Range('B50').Value = 'WOONetwork'
$ws.Range('C50').Value = 'https://coinranking.com/coin/k-J3YwacF+woonetwork-woo'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.2510'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -0.11%  '

$ws.Range('B51').Value = 'Elrond'
$ws.Range('C51').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '34.86'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -0.52%  '
